$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Convert the Mobile numbers in column A (rows 2-16) from text to real numbers ---
$existingMobiles = @{
    2  = 84383321820
    3  = 84396360872
    4  = 84315041690
    5  = 84305303818
    6  = 84345927948
    7  = 84378408033
    8  = 84335987003
    9  = 84344078216
    10 = 84375389598
    11 = 84328861910
    12 = 84377861708
    13 = 84385551213
    14 = 84362423607
    15 = 84354222309
    16 = 84368046503
}

foreach ($r in $existingMobiles.Keys) {
    $ws.Cells.Item($r, 1).Value = $existingMobiles[$r]
}

# --- 2) Append the new registered accounts as rows 17-26 ---
$newAccounts = @(
    @("84377215915", "BUI KHOA TRUONG KINH",   "buikinh537",   "hzxkznO36_", "1954-04-02"),
    @("84307965136", "TRAN THINH HA",          "tranha728",    "dbohtxL72*", "1958-11-27"),
    @("84398006045", "PHAN NGUYEN TRI TAN",    "phantan552",   "buxhrvR32*", "1999-09-21"),
    @("84303354145", "VO THANH HAI DOAN",      "vodoan137",    "egxuvbY22*", "1985-09-06"),
    @("84374198711", "TRAN BINH THINH KHIEU",  "trankhieu660", "jpebzeA56$", "1979-11-02"),
    @("84318857412", "VU THANG TUNG NHAN",     "vunhan216",    "fkjyohJ73@", "1980-09-09"),
    @("84307033964", "BUI CHUNG HUNG",         "buihung546",   "agybobF56$", "1980-11-24"),
    @("84331765380", "VU SAI TRUNG",           "vutrung943",   "yjmmjuE02*", "1997-02-10"),
    @("84317185002", "DANG TAT DINH",          "dangdinh270",  "rxxfjqU77_", "1954-10-02"),
    @("84316963466", "VU VAN DUC THI",         "vuthi211",     "zhvgzmA87@", "1958-12-14")
)

$startRow = 17
$endRow = $startRow + $newAccounts.Length - 1

# Keep the phone numbers (col A) and birthdays (col E) as plain text (matching the
# source data) instead of letting Excel auto-convert them to numbers/dates, so
# format the destination range as text before writing the values into it.
$ws.Range("A$startRow`:E$endRow").NumberFormat = "@"

$r = $startRow
foreach ($account in $newAccounts) {
    $ws.Cells.Item($r, 1).Value = $account[0]
    $ws.Cells.Item($r, 2).Value = $account[1]
    $ws.Cells.Item($r, 3).Value = $account[2]
    $ws.Cells.Item($r, 4).Value = $account[3]
    $ws.Cells.Item($r, 5).Value = $account[4]
    $r++
}
